$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 26,11
$arr[0,0] = 4175
$arr[0,1] = 0
$arr[0,2] = 294
$arr[0,3] = 619
$arr[0,4] = 15
$arr[0,5] = 3
$arr[0,6] = 1015
$arr[0,7] = 26791
$arr[0,8] = 7
$arr[0,9] = 306
$arr[0,10] = 18527
$arr[1,0] = 821
$arr[1,1] = 0
$arr[1,2] = 522
$arr[1,3] = 523
$arr[1,4] = 0
$arr[1,5] = 1
$arr[1,6] = 0
$arr[1,7] = 9093
$arr[1,8] = 11
$arr[1,9] = 0
$arr[1,10] = 0
$arr[2,0] = 256
$arr[2,1] = 0
$arr[2,2] = 0
$arr[2,3] = 0
$arr[2,4] = 0
$arr[2,5] = 0
$arr[2,6] = 0
$arr[2,7] = 0
$arr[2,8] = 0
$arr[2,9] = 0
$arr[2,10] = 0
$arr[3,0] = 803
$arr[3,1] = 0
$arr[3,2] = 212
$arr[3,3] = 214
$arr[3,4] = 0
$arr[3,5] = 2
$arr[3,6] = 0
$arr[3,7] = 8915
$arr[3,8] = 8
$arr[3,9] = 0
$arr[3,10] = 0
$arr[4,0] = 257
$arr[4,1] = 0
$arr[4,2] = 13
$arr[4,3] = 13
$arr[4,4] = 0
$arr[4,5] = 0
$arr[4,6] = 0
$arr[4,7] = 6154
$arr[4,8] = 1
$arr[4,9] = 0
$arr[4,10] = 0
$arr[5,0] = 254
$arr[5,1] = 0
$arr[5,2] = 6
$arr[5,3] = 7
$arr[5,4] = 1
$arr[5,5] = 0
$arr[5,6] = 1429
$arr[5,7] = 1667
$arr[5,8] = 0
$arr[5,9] = 0
$arr[5,10] = 0
$arr[6,0] = 772
$arr[6,1] = 0
$arr[6,2] = 266
$arr[6,3] = 271
$arr[6,4] = 2
$arr[6,5] = 3
$arr[6,6] = 75
$arr[6,7] = 8633
$arr[6,8] = 0
$arr[6,9] = 0
$arr[6,10] = 0
$arr[7,0] = 774
$arr[7,1] = 0
$arr[7,2] = 28
$arr[7,3] = 28
$arr[7,4] = 0
$arr[7,5] = 0
$arr[7,6] = 0
$arr[7,7] = 9643
$arr[7,8] = 0
$arr[7,9] = 0
$arr[7,10] = 0
$arr[8,0] = 738
$arr[8,1] = 2
$arr[8,2] = 48
$arr[8,3] = 55
$arr[8,4] = 5
$arr[8,5] = 0
$arr[8,6] = 1679
$arr[8,7] = 18111
$arr[8,8] = 1
$arr[8,9] = 0
$arr[8,10] = 0
$arr[9,0] = 782
$arr[9,1] = 0
$arr[9,2] = 214
$arr[9,3] = 227
$arr[9,4] = 12
$arr[9,5] = 2
$arr[9,6] = 2568
$arr[9,7] = 48087
$arr[9,8] = 8
$arr[9,9] = 0
$arr[9,10] = 0
$arr[10,0] = 761
$arr[10,1] = 0
$arr[10,2] = 42
$arr[10,3] = 44
$arr[10,4] = 3
$arr[10,5] = 0
$arr[10,6] = 30
$arr[10,7] = 27056
$arr[10,8] = 0
$arr[10,9] = 0
$arr[10,10] = 0
$arr[11,0] = 544
$arr[11,1] = 0
$arr[11,2] = 232
$arr[11,3] = 248
$arr[11,4] = 16
$arr[11,5] = 0
$arr[11,6] = 2183
$arr[11,7] = 4839
$arr[11,8] = 5
$arr[11,9] = 0
$arr[11,10] = 0
$arr[12,0] = 804
$arr[12,1] = 0
$arr[12,2] = 185
$arr[12,3] = 185
$arr[12,4] = 0
$arr[12,5] = 0
$arr[12,6] = 0
$arr[12,7] = 9261
$arr[12,8] = 5
$arr[12,9] = 0
$arr[12,10] = 0
$arr[13,0] = 255
$arr[13,1] = 0
$arr[13,2] = 7
$arr[13,3] = 7
$arr[13,4] = 0
$arr[13,5] = 0
$arr[13,6] = 0
$arr[13,7] = 5714
$arr[13,8] = 0
$arr[13,9] = 0
$arr[13,10] = 0
$arr[14,0] = 289
$arr[14,1] = 0
$arr[14,2] = 338
$arr[14,3] = 432
$arr[14,4] = 5
$arr[14,5] = 2
$arr[14,6] = 421
$arr[14,7] = 13843
$arr[14,8] = 2
$arr[14,9] = 87
$arr[14,10] = 2358
$arr[15,0] = 47
$arr[15,1] = 0
$arr[15,2] = 7
$arr[15,3] = 8
$arr[15,4] = 1
$arr[15,5] = 0
$arr[15,6] = 1250
$arr[15,7] = 4286
$arr[15,8] = 0
$arr[15,9] = 0
$arr[15,10] = 0
$arr[16,0] = 215
$arr[16,1] = 0
$arr[16,2] = 320
$arr[16,3] = 362
$arr[16,4] = 4
$arr[16,5] = 0
$arr[16,6] = 110
$arr[16,7] = 6006
$arr[16,8] = 2
$arr[16,9] = 39
$arr[16,10] = 1077
$arr[17,0] = 1511
$arr[17,1] = 0
$arr[17,2] = 124
$arr[17,3] = 125
$arr[17,4] = 1
$arr[17,5] = 0
$arr[17,6] = 95
$arr[17,7] = 12538
$arr[17,8] = 0
$arr[17,9] = 0
$arr[17,10] = 0
$arr[18,0] = 1374
$arr[18,1] = 0
$arr[18,2] = 98
$arr[18,3] = 131
$arr[18,4] = 33
$arr[18,5] = 0
$arr[18,6] = 5377
$arr[18,7] = 6717
$arr[18,8] = 2
$arr[18,9] = 0
$arr[18,10] = 0
$arr[19,0] = 416
$arr[19,1] = 0
$arr[19,2] = 52
$arr[19,3] = 51
$arr[19,4] = 0
$arr[19,5] = 0
$arr[19,6] = 0
$arr[19,7] = 30784
$arr[19,8] = 0
$arr[19,9] = 0
$arr[19,10] = 0
$arr[20,0] = 177
$arr[20,1] = 0
$arr[20,2] = 4
$arr[20,3] = 4
$arr[20,4] = 0
$arr[20,5] = 0
$arr[20,6] = 0
$arr[20,7] = 0
$arr[20,8] = 0
$arr[20,9] = 0
$arr[20,10] = 0
$arr[21,0] = 813
$arr[21,1] = 0
$arr[21,2] = 6
$arr[21,3] = 6
$arr[21,4] = 0
$arr[21,5] = 0
$arr[21,6] = 0
$arr[21,7] = 6667
$arr[21,8] = 0
$arr[21,9] = 0
$arr[21,10] = 0
$arr[22,0] = 938
$arr[22,1] = 0
$arr[22,2] = 13
$arr[22,3] = 14
$arr[22,4] = 1
$arr[22,5] = 0
$arr[22,6] = 1429
$arr[22,7] = 15476
$arr[22,8] = 0
$arr[22,9] = 0
$arr[22,10] = 0
$arr[23,0] = 2
$arr[23,1] = 1
$arr[23,2] = 29
$arr[23,3] = 54
$arr[23,4] = 3
$arr[23,5] = 0
$arr[23,6] = 556
$arr[23,7] = 1379
$arr[23,8] = 0
$arr[23,9] = 20
$arr[23,10] = 3704
$arr[24,0] = 118
$arr[24,1] = 0
$arr[24,2] = 21
$arr[24,3] = 26
$arr[24,4] = 4
$arr[24,5] = 0
$arr[24,6] = 1538
$arr[24,7] = 4286
$arr[24,8] = 0
$arr[24,9] = 0
$arr[24,10] = 0
$arr[25,0] = 173
$arr[25,1] = 0
$arr[25,2] = 71
$arr[25,3] = 76
$arr[25,4] = 5
$arr[25,5] = 0
$arr[25,6] = 658
$arr[25,7] = 1690
$arr[25,8] = 1
$arr[25,9] = 0
$arr[25,10] = 0
$ws.Range("B2:L27").Value = $arr
Write-Output "done"
